$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 4).Value = 2.499999999999988
$ws.Cells.Item(2, 8).Value = 200000
$ws.Cells.Item(3, 4).Value = 2.499999999999988
$ws.Cells.Item(3, 8).Value = 200000
$ws.Cells.Item(4, 4).Value = 2.499999999999988
$ws.Cells.Item(4, 8).Value = 200000
$ws.Cells.Item(5, 4).Value = 2.499999999999988
$ws.Cells.Item(5, 8).Value = 200000
$ws.Cells.Item(6, 4).Value = 2.499999999999988
$ws.Cells.Item(6, 8).Value = 200000
$ws.Cells.Item(7, 4).Value = 9.99999999999995
$ws.Cells.Item(7, 8).Value = 200000
$ws.Cells.Item(8, 4).Value = 9.99999999999995
$ws.Cells.Item(8, 8).Value = 200000
$ws.Cells.Item(9, 4).Value = 2.499999999999988
$ws.Cells.Item(9, 8).Value = 200000
$ws.Cells.Item(10, 4).Value = 14.99999999999993
$ws.Cells.Item(10, 8).Value = 200000

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 4).Value = 2.499999999999988
$ws.Cells.Item(2, 8).Value = 200000
$ws.Cells.Item(3, 4).Value = 2.499999999999988
$ws.Cells.Item(3, 8).Value = 200000
$ws.Cells.Item(4, 4).Value = 2.499999999999988
$ws.Cells.Item(4, 8).Value = 200000
$ws.Cells.Item(5, 4).Value = 0.4999999999999975
$ws.Cells.Item(5, 8).Value = 200000
$ws.Cells.Item(6, 4).Value = 22.49999999999989
$ws.Cells.Item(6, 8).Value = 200000
$ws.Cells.Item(7, 4).Value = 4.499999999999978
$ws.Cells.Item(7, 8).Value = 200000
$ws.Cells.Item(8, 4).Value = 14.99999999999993
$ws.Cells.Item(8, 8).Value = 200000

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 4).Value = 0.4999999999999975
$ws.Cells.Item(2, 8).Value = 200000
$ws.Cells.Item(3, 4).Value = 2.499999999999988
$ws.Cells.Item(3, 8).Value = 200000
$ws.Cells.Item(4, 4).Value = 2.499999999999988
$ws.Cells.Item(4, 8).Value = 200000
$ws.Cells.Item(5, 4).Value = 2.499999999999988
$ws.Cells.Item(5, 8).Value = 200000
$ws.Cells.Item(6, 4).Value = 2.499999999999988
$ws.Cells.Item(6, 8).Value = 200000
$ws.Cells.Item(7, 4).Value = 0.4999999999999975
$ws.Cells.Item(7, 8).Value = 200000
$ws.Cells.Item(8, 4).Value = 0.4999999999999975
$ws.Cells.Item(8, 8).Value = 200000
$ws.Cells.Item(9, 4).Value = 0.4999999999999975
$ws.Cells.Item(9, 8).Value = 200000
$ws.Cells.Item(10, 4).Value = 0.4999999999999975
$ws.Cells.Item(10, 8).Value = 200000
$ws.Cells.Item(11, 4).Value = 0.4999999999999975
$ws.Cells.Item(11, 8).Value = 200000
$ws.Cells.Item(12, 4).Value = 0.4999999999999975
$ws.Cells.Item(12, 8).Value = 200000
$ws.Cells.Item(13, 4).Value = 0.4999999999999975
$ws.Cells.Item(13, 8).Value = 200000
$ws.Cells.Item(14, 4).Value = 0.4999999999999975
$ws.Cells.Item(14, 8).Value = 200000
$ws.Cells.Item(15, 4).Value = 35.49999999999983
$ws.Cells.Item(15, 8).Value = 200000
